$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @("特变电工", "特变电工", "金风科技")
    3  = @("华胜天成", "华胜天成", "华胜天成")
    4  = @("长电科技", "金风科技", "海格通信")
    5  = @("保变电气", "岩山科技", "岩山科技")
    6  = @("金风科技", "保变电气", "航天发展")
    7  = @("通富微电", "航天发展", "利欧股份")
    8  = @("岩山科技", "长电科技", "特变电工")
    9  = @("蓝色光标", "中国西电", "蓝色光标")
    10 = @("航天发展", "三花智控", "兆易创新")
    11 = @("利欧股份", "利欧股份", "美年健康")
    12 = @("五洲新春", "蓝色光标", "平潭发展")
    13 = @("海格通信", "通富微电", "长电科技")
    14 = @("兆易创新", "兆易创新", "通富微电")
    15 = @("三花智控", "五洲新春", "三维通信")
    16 = @("中国西电", "海格通信", "盈新发展")
    17 = @("中国卫星", "山子高科", "国晟科技")
    18 = @("雷科防务", "航天电子", "雷科防务")
    19 = @("航天电子", "雷科防务", "神剑股份")
    20 = @("三变科技", "中国卫星", "华夏幸福")
    21 = @("圣晖集成", "三变科技", "三花智控")
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
}
